$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the date columns (F,H,I,J,L,N,O,P,Q,S,U,V,X,Y,Z) as Text so the
# new values are stored as shared strings instead of date serials.
$ws.Range("F1:F4").NumberFormat = "@"
$ws.Range("H1:H4").NumberFormat = "@"
$ws.Range("I1:I4").NumberFormat = "@"
$ws.Range("J1:J4").NumberFormat = "@"
$ws.Range("L1:L4").NumberFormat = "@"
$ws.Range("N1:N4").NumberFormat = "@"
$ws.Range("O1:O4").NumberFormat = "@"
$ws.Range("P1:P4").NumberFormat = "@"
$ws.Range("Q1:Q4").NumberFormat = "@"
$ws.Range("S1:S4").NumberFormat = "@"
$ws.Range("U1:U4").NumberFormat = "@"
$ws.Range("V1:V4").NumberFormat = "@"
$ws.Range("X1:X4").NumberFormat = "@"
$ws.Range("Y1:Y4").NumberFormat = "@"
$ws.Range("Z1:Z4").NumberFormat = "@"

# Row 2 - replace serial dates with text "DD-Mon-YYYY" values (PHA Dt for
# this patient was updated, hence 01-Jan-2024 in column N).
$ws.Range("F2").Value = "01-Jan-2000"
$ws.Range("H2").Value = "02-Jan-2000"
$ws.Range("I2").Value = "03-Jan-2000"
$ws.Range("J2").Value = "04-Jan-2000"
$ws.Range("L2").Value = "05-Jan-2000"
$ws.Range("N2").Value = "01-Jan-2024"
$ws.Range("O2").Value = "07-Jan-2000"
$ws.Range("P2").Value = "08-Jan-2000"
$ws.Range("Q2").Value = "09-Jan-2000"
$ws.Range("S2").Value = "10-Jan-2000"
$ws.Range("U2").Value = "11-Jan-2000"
$ws.Range("V2").Value = "12-Jan-2000"
$ws.Range("X2").Value = "13-Jan-2000"
$ws.Range("Y2").Value = "14-Jan-2000"
$ws.Range("Z2").Value = "15-Jan-2000"

# Row 3
$ws.Range("F3").Value = "01-Jan-2000"
$ws.Range("H3").Value = "02-Jan-2000"
$ws.Range("I3").Value = "03-Jan-2000"
$ws.Range("J3").Value = "04-Jan-2000"
$ws.Range("L3").Value = "05-Jan-2000"
$ws.Range("N3").Value = "06-Jan-2000"
$ws.Range("O3").Value = "07-Jan-2000"
$ws.Range("P3").Value = "08-Jan-2000"
$ws.Range("Q3").Value = "09-Jan-2000"
$ws.Range("S3").Value = "10-Jan-2000"
$ws.Range("U3").Value = "11-Jan-2000"
$ws.Range("V3").Value = "12-Jan-2000"
$ws.Range("X3").Value = "13-Jan-2000"
$ws.Range("Y3").Value = "14-Jan-2000"
$ws.Range("Z3").Value = "15-Jan-2000"

# Row 4 - N4 stays blank (PHA Dt missing for this patient), so just clear
# its old (empty, date-formatted) cell rather than giving it a value.
$ws.Range("F4").Value = "01-Jan-2000"
$ws.Range("H4").Value = "02-Jan-2000"
$ws.Range("I4").Value = "03-Jan-2000"
$ws.Range("J4").Value = "04-Jan-2000"
$ws.Range("L4").Value = "05-Jan-2000"
$ws.Range("N4").Clear()
$ws.Range("O4").Value = "07-Jan-2000"
$ws.Range("P4").Value = "08-Jan-2000"
$ws.Range("Q4").Value = "09-Jan-2000"
$ws.Range("S4").Value = "10-Jan-2000"
$ws.Range("U4").Value = "11-Jan-2000"
$ws.Range("V4").Value = "12-Jan-2000"
$ws.Range("X4").Value = "13-Jan-2000"
$ws.Range("Y4").Value = "14-Jan-2000"
$ws.Range("Z4").Value = "15-Jan-2000"

# Leave the selection where the editing session ended.
$ws.Range("K4").Select()
